$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates: force Text number format so the
# numeric-looking strings are preserved as text, matching the
# source workbook where these cells are stored as inline/shared strings.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D22","D23","D25","D26","D40","D42","D43","D44","D45","D46","D47","D48","D49")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "246.76"
$ws.Range("D3").Value = "21.81"
$ws.Range("D4").Value = "5.464"
$ws.Range("D5").Value = "0.05661"
$ws.Range("D6").Value = "3.375"
$ws.Range("D7").Value = "0.8017"
$ws.Range("D8").Value = "1.036"
$ws.Range("D9").Value = "0.01158"
$ws.Range("D10").Value = "0.1442"
$ws.Range("D11").Value = "0.07252"
$ws.Range("D12").Value = "0.03162"
$ws.Range("D13").Value = "0.02941"
$ws.Range("D14").Value = "0.09283"
$ws.Range("D15").Value = "0.001645"
$ws.Range("D16").Value = "3.208"
$ws.Range("D17").Value = "0.04703"
$ws.Range("D18").Value = "0.006437"
$ws.Range("D19").Value = "0.005038"
$ws.Range("D20").Value = "0.001046"
$ws.Range("D22").Value = "0.0003202"
$ws.Range("D23").Value = "3.817"
$ws.Range("D25").Value = "2.089"
$ws.Range("D26").Value = "0.3330"
$ws.Range("D40").Value = "0.04087"
$ws.Range("D42").Value = "0.003502"
$ws.Range("D43").Value = "0.1040"
$ws.Range("D44").Value = "0.008051"
$ws.Range("D45").Value = "0.00005851"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").Value = "0.6829"
$ws.Range("D48").Value = "0.01032"
$ws.Range("D49").Value = "0.00002101"

# --- Coin / Link / Volume(1h) text updates (rows reshuffled as the
# scraped ranking list moved by one position, plus a couple of swaps).
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E9").Value = "8OneONEBestin24h"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"
